$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the used range so we know how many rows to touch.
$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

# 1) Update column C ("Förändrad") from 45184 -> 45186 for every data row.
for ($i = 2; $i -le $lastRow; $i++) {
    $cCell = $ws.Cells.Item($i, 3)
    $cVal = $cCell.Value2
    if ($cVal -eq 45184) {
        $cCell.Value = 45186
    }
}

# 2) Add the case number (column A) as a second HYPERLINK() argument for
#    every hyperlink-formula cell (columns S, T, U, V, W, X, Y).
$hyperlinkCols = @("S","T","U","V","W","X","Y")
for ($i = 2; $i -le $lastRow; $i++) {
    $aVal = $ws.Cells.Item($i, 1).Value2
    if ($aVal -eq $null -or $aVal -eq "") {
        continue
    }
    foreach ($col in $hyperlinkCols) {
        $cell = $ws.Range($col + $i)
        $f = $cell.Formula
        if ($f -ne $null -and $f.Length -gt 0 -and $f.StartsWith("=HYPERLINK(") -and -not $f.Contains(",")) {
            $newFormula = $f.Substring(0, $f.Length - 1) + ', "' + $aVal + '")'
            $cell.Formula = $newFormula
        }
    }
}
